$wb = $excel.ActiveWorkbook

# --- ColumnHeadersNcp sheet: fix temperature definition text ---
$wsNcp = $wb.Worksheets.Item("ColumnHeadersNcp")
$wsNcp.Range("B7").Value = "Underway thermosalinograph temperature in degrees Celsius. URI http://vocab.nerc.ac.uk/collection/P01/current/TEMPSZ01/"
$wsNcp.Range("B7").Select()

# --- Personnel sheet: add Kate Morkeski as metadata Provider ---
$wsPer = $wb.Worksheets.Item("Personnel")
$wsPer.Range("A9").Value = "Kate"
$wsPer.Range("C9").Value = "Morkeski"
$wsPer.Range("D9").Value = "Northeast U.S. Shelf LTER"
$wsPer.Range("E9").Value = "kmorkeski@whoi.edu"
$wsPer.Range("F9").Value = "0000-0002-2903-5851"
$wsPer.Range("G9").Value = "metadata Provider"
$wsPer.Range("H9").Value = "Northeast U.S. Shelf LTER"
$wsPer.Range("I9").Value = "NSF"
$wsPer.Range("J9").Value = "OCE-2322676"
$wsPer.Range("A9:J9").Select()

$wb.Save()
